$wb = $excel.ActiveWorkbook

# --- 展览 (sheet 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 203
$ws.Range("F4").Value = 369
$ws.Range("F5").Value = 1634
$ws.Range("F6").Value = 810
$ws.Range("F7").Value = 703
$ws.Range("F8").Value = 1291
$ws.Range("F9").Value = 2641
$ws.Range("F10").Value = 1351
$ws.Range("C11").Value = "上海·原神×星穹铁道ONLY 2.0"
$ws.Range("D11").Value = "吴中路1588号上海爱琴海购物中心F4 竞梦元宇宙"
$ws.Range("F11").Value = 2052
$ws.Range("G11").Value = 68
$ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=88273"
$ws.Range("I11").Value = "//i0.hdslb.com/bfs/openplatform/202407/wcRAjLG11721196920833.png"
$ws.Range("C12").Value = "上海·漫游L+动漫游戏嘉年华（免费展）"
$ws.Range("D12").Value = "申长路869号 上海龙湖虹桥天街"
$ws.Range("F12").Value = 828
$ws.Range("G12").Value = 20
$ws.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=88134"
$ws.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202406/adaE6Z6f1719454819535.jpeg"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "2024-07-20"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "上海·第九届Redamancy动漫游戏嘉年华"
$ws.Range("D13").Value = "中山北路3300号4楼 上海环球港"
$ws.Range("E13").Value = "2024.07.20 10:00-07.21 17:00"
$ws.Range("F13").Value = 2323
$ws.Range("G13").Value = 60
$ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=84637"
$ws.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202404/hWLkXqwM1713194236349.png"
$ws.Range("F14").Value = 721
$ws.Range("F15").Value = 6567
$ws.Range("F17").Value = 515
$ws.Range("F18").Value = 1240
$ws.Range("F20").Value = 1494
$ws.Range("F21").Value = 1333
$ws.Range("F22").Value = 1204
$ws.Range("F23").Value = 98
$ws.Range("F24").Value = 2394
$ws.Range("F25").Value = 1106
$ws.Range("F26").Value = 1009
$ws.Range("F27").Value = 749
$ws.Range("F28").Value = 1117
$ws.Range("F29").Value = 248
$ws.Range("F30").Value = 5328
$ws.Range("F32").Value = 898
$ws.Range("F33").Value = 1258
$ws.Range("F35").Value = 3727
$ws.Range("F37").Value = 1693
$ws.Range("F38").Value = 166
$ws.Range("F39").Value = 963
$ws.Range("F40").Value = 1052
$ws.Range("F41").Value = 400
$ws.Range("F42").Value = 1772
$ws.Range("F43").Value = 906
$ws.Range("F44").Value = 1047
$ws.Range("F45").Value = 509
$ws.Range("F46").Value = 514
$ws.Range("F47").Value = 13
$ws.Range("F48").Value = 61
$ws.Range("F49").Value = 84
# --- 演出 (sheet 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F11").Value = 393
$ws.Range("F13").Value = 141
$ws.Range("F21").Value = 593
$ws.Range("F22").Value = 246
$ws.Range("F26").Value = 83
$ws.Range("F27").Value = 83
$ws.Range("F30").Value = 304
$ws.Range("F31").Value = 46
$ws.Range("F35").Value = 47
$ws.Range("F37").Value = 112
$ws.Range("F39").Value = 201
# --- 本地生活 (sheet 3) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 3298
$ws.Range("F5").Value = 398
$ws.Range("F7").Value = 1460
$ws.Range("F8").Value = 766
$ws.Range("F9").Value = 388
$ws.Range("F10").Value = 2803
$ws.Range("F11").Value = 298
$ws.Range("F12").Value = 545
$ws.Range("F13").Value = 626
$ws.Range("F14").Value = 1174
# --- 全部类型 (sheet 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 398
$ws.Range("F3").Value = 766
$ws.Range("F4").Value = 369
$ws.Range("F5").Value = 388
$ws.Range("F6").Value = 2803
$ws.Range("F7").Value = 1635
$ws.Range("F10").Value = 703
$ws.Range("F11").Value = 1291
$ws.Range("F12").Value = 2641
$ws.Range("F13").Value = 1351
$ws.Range("C14").Value = "上海·原神×星穹铁道ONLY 2.0"
$ws.Range("D14").Value = "吴中路1588号上海爱琴海购物中心F4 竞梦元宇宙"
$ws.Range("F14").Value = 2052
$ws.Range("G14").Value = 68
$ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=88273"
$ws.Range("I14").Value = "//i0.hdslb.com/bfs/openplatform/202407/wcRAjLG11721196920833.png"
$ws.Range("C15").Value = "上海·漫游L+动漫游戏嘉年华（免费展）"
$ws.Range("D15").Value = "申长路869号 上海龙湖虹桥天街"
$ws.Range("F15").Value = 828
$ws.Range("G15").Value = 20
$ws.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=88134"
$ws.Range("I15").Value = "//i1.hdslb.com/bfs/openplatform/202406/adaE6Z6f1719454819535.jpeg"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "2024-07-20"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "上海·第九届Redamancy动漫游戏嘉年华"
$ws.Range("D16").Value = "中山北路3300号4楼 上海环球港"
$ws.Range("E16").Value = "2024.07.20 10:00-07.21 17:00"
$ws.Range("F16").Value = 2323
$ws.Range("G16").Value = 60
$ws.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=84637"
$ws.Range("I16").Value = "//i1.hdslb.com/bfs/openplatform/202404/hWLkXqwM1713194236349.png"
$ws.Range("F17").Value = 721
$ws.Range("F18").Value = 6568
$ws.Range("F19").Value = 545
$ws.Range("F20").Value = 515
$ws.Range("F21").Value = 1240
$ws.Range("F22").Value = 626
$ws.Range("F24").Value = 1333
$ws.Range("F25").Value = 98
$ws.Range("F26").Value = 2394
$ws.Range("F27").Value = 246
$ws.Range("F28").Value = 83
$ws.Range("F29").Value = 1106
$ws.Range("F30").Value = 1117
$ws.Range("F31").Value = 248
$ws.Range("F32").Value = 5328
$ws.Range("F34").Value = 901
$ws.Range("F35").Value = 1258
$ws.Range("F36").Value = 3727
$ws.Range("F38").Value = 304
$ws.Range("F39").Value = 1693
$ws.Range("F40").Value = 166
$ws.Range("F41").Value = 963
$ws.Range("F42").Value = 1772
$ws.Range("F43").Value = 906
$ws.Range("F44").Value = 1047
$ws.Range("F45").Value = 509
$ws.Range("F46").Value = 514
$ws.Range("F47").Value = 201
$ws.Range("F48").Value = 201
$ws.Range("F49").Value = 13
$ws.Range("F50").Value = 61
$ws.Range("F51").Value = 84
